$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "gemini-1.5-pro"
$ws.Range("C9").Value = "-0.1 ± 0.55"
$ws.Range("D9").Value = "0.44 ± 0.79"
$ws.Range("J9").Value = "0.82 ± 0.1"
$ws.Range("K9").Value = "0.85 ± 0.11"
$ws.Range("L9").Value = "0.84 ± 0.11"
$ws.Range("M9").Value = "0.9 ± 0.12"
$ws.Range("N9").Value = "0.98 ± 0.12"
$ws.Range("P9").Value = "0.48 ± 0.14"
$ws.Range("Q9").Value = "5.21 ± 5.4"
$ws.Range("R9").Value = "0.015 ± 0.00"
$ws.Range("S9").Value = "0.92 ± 0.12"
$ws.Range("T9").Value = "0.97 ± 0.15"
$ws.Range("U9").Value = "2.96 ± 1.27"
$ws.Range("V9").Value = "0.54 ± 0.41"
$ws.Range("W9").Value = "0.92 ± 0.12"
$ws.Range("X9").Value = "1.26 ± 0.31"
